# Update calibration cost values for rows 100-107 and 114-115, columns J:AS.
# Each row holds a single repeated value across columns J through AS (36 cols).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 37774.29676
    101 = 283528.6032
    102 = 15524.95214
    103 = 11633.79687
    104 = 29057.15135
    105 = 1054.207668
    106 = 1707.744735
    107 = 21786.22172
    114 = 27.6712064
    115 = 13123.17814
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $rng = $ws.Range("J$row`:AS$row")
    $rng.Value = $value
}
